$wb = $excel.ActiveWorkbook

# --- Sheet "Concussions" (sheet1): add 2019 row (row 10) ---
$ws1 = $wb.Worksheets.Item("Concussions")
$ws1.Range("A10").Value = 2019
$ws1.Range("B10").Value = 30
$ws1.Range("C10").Value = 49
$ws1.Range("D10").Formula = "=B10+C10"
$ws1.Range("E10").Value = 9
$ws1.Range("F10").Value = 136
$ws1.Range("G10").Formula = "=F10+E10"
$ws1.Range("H10").Formula = "=E10+B10"
$ws1.Range("I10").Formula = "=F10+C10"
$ws1.Range("J10").Formula = "=I10+H10"
$ws1.Activate()
$ws1.Range("A10").Select()

# --- Sheet "ACL Tears" (sheet2): add 2019 row (row 10) ---
$ws2 = $wb.Worksheets.Item("ACL Tears")
$ws2.Range("A10").Value = 2019
$ws2.Range("B10").Value = 7
$ws2.Range("C10").Value = 10
$ws2.Range("D10").Formula = "=B10+C10"
$ws2.Range("E10").Value = 7
$ws2.Range("F10").Value = 23
$ws2.Range("G10").Formula = "=F10+E10"
$ws2.Range("H10").Formula = "=E10+B10"
$ws2.Range("I10").Formula = "=F10+C10"
$ws2.Range("J10").Formula = "=I10+H10"
$ws2.Activate()
$ws2.Range("A10").Select()

# --- Sheet "MCL Tears" (sheet3): add 2019 row (row 10) ---
$ws3 = $wb.Worksheets.Item("MCL Tears")
$ws3.Range("A10").Value = 2019
$ws3.Range("B10").Value = 7
$ws3.Range("C10").Value = 23
$ws3.Range("D10").Formula = "=B10+C10"
$ws3.Range("E10").Value = 3
$ws3.Range("F10").Value = 76
$ws3.Range("G10").Formula = "=F10+E10"
$ws3.Range("H10").Formula = "=E10+B10"
$ws3.Range("I10").Formula = "=F10+C10"
$ws3.Range("J10").Formula = "=I10+H10"
$ws3.Activate()
$ws3.Range("G10").Select()
